# ----------------------------------------------------------------------------
# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp (A1)
# - A handful of countries leap-frog their neighbours in the ranking, so their
#   name swaps with the row above/below while the stats columns are refreshed
# - Refresh case/death/recovery counters (cols B:H) for the rows with new data
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 26 de Agosto de 2020 a las 19:22"

# Country name re-shuffle (A column) --------------------------------------------
$nameUpdates = @{
    55 = "Etiopia"
    56 = "Ghana"
    93 = "Grecia"
    94 = "Guinea"
    105 = "Namibia"
    106 = "Zimbabue"
    119 = "Mozambique"
    120 = "Cabo Verde"
    121 = "Ruanda"
    122 = "Eslovaquia"
    143 = "Jordania"
    144 = "Malta"
    145 = "Jamaica"
}
foreach ($row in $nameUpdates.Keys) {
    $ws.Cells.Item($row, 1).Value = $nameUpdates[$row]
}

# Updated statistics (B:Casos totales, C:Nuevos casos, D:Casos activos,
# E:Recuperados, F:Casos criticos, G:Muertes hoy, H:Muertes) ---------------------
$statUpdates = @(
    @(4,2,5970799),
    @(4,3,15071),
    @(4,4,3257962),
    @(4,5,2529914),
    @(4,7,519),
    @(4,8,182923),
    @(5,2,3683224),
    @(5,3,9048),
    @(5,5,717865),
    @(5,7,298),
    @(5,8,116964),
    @(6,2,3288693),
    @(6,3,56939),
    @(6,4,2512544),
    @(6,5,715852),
    @(6,7,685),
    @(6,8,60297),
    @(12,2,426818),
    @(12,3,3594),
    @(12,7,47),
    @(12,8,28971),
    @(13,2,402365),
    @(13,3,1380),
    @(13,4,376268),
    @(13,5,15107),
    @(13,7,32),
    @(13,8,10990),
    @(16,2,328846),
    @(16,3,1048),
    @(16,7,16),
    @(16,8,41465),
    @(21,2,262507),
    @(21,3,1313),
    @(21,4,239797),
    @(21,5,16527),
    @(21,7,20),
    @(21,8,6183),
    @(24,2,215784),
    @(24,3,3837),
    @(24,4,157215),
    @(24,5,51901),
    @(24,7,72),
    @(24,8,6668),
    @(32,2,108054),
    @(32,3,1594),
    @(32,4,86450),
    @(32,5,20729),
    @(32,7,16),
    @(32,8,875),
    @(35,2,92557),
    @(35,3,340),
    @(35,4,63478),
    @(35,5,27466),
    @(35,7,28),
    @(35,8,1613),
    @(55,2,45221),
    @(55,3,1533),
    @(55,4,16311),
    @(55,5,28185),
    @(55,7,16),
    @(55,8,725),
    @(56,2,43717),
    @(56,4,41843),
    @(56,5,1604),
    @(56,8,270),
    @(70,2,28363),
    @(70,3,162),
    @(70,5,3222),
    @(74,2,22790),
    @(74,3,242),
    @(74,4,16939),
    @(74,5,5433),
    @(74,7,2),
    @(74,8,418),
    @(89,2,11376),
    @(89,3,91),
    @(89,4,10693),
    @(89,5,401),
    @(93,2,9280),
    @(93,3,293),
    @(93,4,3804),
    @(93,5,5228),
    @(93,7,5),
    @(93,8,248),
    @(94,2,9167),
    @(94,3,39),
    @(94,4,8150),
    @(94,5,960),
    @(94,8,57),
    @(103,2,7225),
    @(103,3,178),
    @(103,4,4561),
    @(103,5,2636),
    @(105,2,6431),
    @(105,3,271),
    @(105,4,2734),
    @(105,5,3638),
    @(105,7,2),
    @(105,8,59),
    @(106,2,6196),
    @(106,4,4961),
    @(106,5,1069),
    @(106,8,166),
    @(119,2,3590),
    @(119,3,82),
    @(119,4,1927),
    @(119,5,1642),
    @(119,8,21),
    @(120,2,3568),
    @(120,4,2673),
    @(120,5,858),
    @(120,8,37),
    @(121,2,3537),
    @(121,3,0),
    @(121,4,1806),
    @(121,5,1716),
    @(121,8,15),
    @(122,2,3536),
    @(122,3,84),
    @(122,4,2192),
    @(122,5,1311),
    @(122,8,33),
    @(131,2,2708),
    @(131,3,22),
    @(131,4,611),
    @(131,5,2004),
    @(131,7,3),
    @(131,8,93),
    @(133,2,2510),
    @(133,3,3),
    @(133,5,1173),
    @(140,2,2003),
    @(140,3,2),
    @(140,4,1577),
    @(140,5,357),
    @(143,2,1756),
    @(143,3,40),
    @(143,4,1355),
    @(143,5,386),
    @(143,7,1),
    @(143,8,15),
    @(144,2,1751),
    @(144,3,46),
    @(144,4,1077),
    @(144,5,664),
    @(144,7,0),
    @(144,8,10),
    @(145,2,1732),
    @(145,3,120),
    @(145,4,840),
    @(145,5,873),
    @(145,7,3),
    @(145,8,19)
)
foreach ($u in $statUpdates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}
